$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "<tho>"
$ws.Range("C2").Value = 32

$ws.Range("C3").Value = 33

$ws.Range("C4").Value = 35

$ws.Range("C5").Value = 32

$ws.Range("C6").Value = 31

$ws.Range("C7").Value = 32

$ws.Range("B8").Value = "<num>"
$ws.Range("C8").Value = 38

$ws.Range("B9").Value = "<foxtrot>"
$ws.Range("C9").Value = 37

$ws.Range("B10").Value = "<xaul>"
$ws.Range("C10").Value = 34

$ws.Range("B11").Value = "<sie>"
$ws.Range("C11").Value = 36

$ws.Range("C13").Value = 36

$ws.Range("B14").Value = "<november>"
$ws.Range("C14").Value = 36

$ws.Range("C15").Value = 31

$ws.Range("B16").Value = "<numre>"
$ws.Range("C16").Value = 32

$ws.Range("C17").Value = 34

$ws.Range("B18").Value = "<whit>"
$ws.Range("C18").Value = 28
